# Pushes the `indiv` tax/medicare table back to cover 2002-03 .. 1999-00,
# fixing assorted typos, per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("indiv")

# Exact number-format codes reused from the existing sheet so that we reuse
# the same style indexes (s="1" comma-0dp, s="3" comma-3dp) instead of Excel
# minting brand new (duplicate) number formats.
$fmtComma0 = "_-* #,##0_-;\-* #,##0_-;_-* ""-""??_-;_-@_-"
$fmtComma3 = "_-* #,##0.000_-;\-* #,##0.000_-;_-* ""-""??_-;_-@_-"

# Column indices: A=1 fy_year, B=2 sato, C=3 pto, D=4 sapto, E=5 family_status,
# F=6 lower_threshold, G=7 upper_threshold, H=8 taper, I=9 rate,
# J=10 lower_family_threshold, K=11 upper_family_threshold,
# L=12 lower_up_for_each_child

# Each entry: row, fy_year, sato, pto, sapto, family_status, F, H, I, J, K, L
# ($null means "leave the cell blank")
$rows = @(
  @(39, "2002-03", $false, $false, $false, "individual", 15062, 0.2, 0.015, 25417, 25418, 2334),
  @(40, "2002-03", $true,  $false, $true,  "individual", 20000, 0.2, 0.015, 31729, 31730, 2334),
  @(41, "2002-03", $false, $true,  $false, "individual", 17164, 0.2, 0.015, 31729, 31730, 2334),

  @(42, "2001-02", $false, $false, $false, "individual", 14539, 0.2, 0.015, 24534, 24535, 2253),
  @(43, "2001-02", $true,  $false, $true,  "individual", 20000, 0.2, 0.015, 31729, 31730, 2253),
  @(44, "2001-02", $false, $true,  $false, "individual", 16570, 0.2, 0.015, 31729, 31730, 2253),

  @(45, "2000-01", $false, $false, $false, "individual", 13807, 0.2, 0.015, 23300, 23301, 2140),
  @(46, "2000-01", $true,  $false, $true,  "individual", 20000, 0.2, 0.015, 31729, 31730, 2140),
  @(47, "2000-01", $false, $true,  $false, "individual", 15970, 0.2, 0.015, 31729, 31730, 2140),

  @(48, "1999-00", $false, $false, $false, "individual", 13351, 0.2, 0.015, $null, $null, $null),
  @(49, "1999-00", $false, $false, $false, "Spouse without dependant child or student child", 22866, 0.2, 0.015, $null, $null, $null),
  @(50, "1999-00", $false, $false, $false, "Spouse with one dependant child", 24966, 0.2, 0.015, $null, $null, $null),
  @(51, "1999-00", $false, $false, $false, "Spouse with more than one dependant child", 24966, 0.2, 0.015, $null, $null, $null),

  @(52, "1999-00", $false, $false, $false, $null, $null, $null, $null, $null, $null, $null),
  @(53, "1999-00", $true,  $false, $true,  $null, $null, $null, $null, $null, $null, $null),
  @(54, "1999-00", $false, $true,  $false, $null, $null, $null, $null, $null, $null, $null)
)

foreach ($r in $rows) {
    $rowNum = $r[0]

    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]

    if ($null -ne $r[5]) {
        $ws.Cells.Item($rowNum, 5).Value = $r[5]
    }
    if ($null -ne $r[6]) {
        $ws.Cells.Item($rowNum, 6).Value = $r[6]
        $ws.Cells.Item($rowNum, 6).NumberFormat = $fmtComma0
    }
    if ($null -ne $r[7]) {
        $ws.Cells.Item($rowNum, 8).Value = $r[7]
        $ws.Cells.Item($rowNum, 8).NumberFormat = $fmtComma3
    }
    if ($null -ne $r[8]) {
        $ws.Cells.Item($rowNum, 9).Value = $r[8]
    }
    if ($null -ne $r[9]) {
        $ws.Cells.Item($rowNum, 10).Value = $r[9]
        $ws.Cells.Item($rowNum, 10).NumberFormat = $fmtComma0
    }
    if ($null -ne $r[10]) {
        $ws.Cells.Item($rowNum, 11).Value = $r[10]
        $ws.Cells.Item($rowNum, 11).NumberFormat = $fmtComma0
    }
    if ($null -ne $r[11]) {
        $ws.Cells.Item($rowNum, 12).Value = $r[11]
        $ws.Cells.Item($rowNum, 12).NumberFormat = $fmtComma0
    }
}

# G38:G54 becomes one shared formula (matches the existing shared-formula
# pattern already used higher up the sheet, e.g. si="2", si="3").
$ws.Range("G38:G54").Formula = "=ROUND((F38 +1)*H38/(H38-I38), 0)"
$ws.Range("G38:G54").NumberFormat = $fmtComma0

# Freeze the header row and scroll/select so row 38 onward is visible, with
# F47 as the active cell (as in the authored workbook).
$app = $ws.Application
$win = $app.ActiveWindow
$ws.Range("A2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 38
$ws.Range("F47").Select()
